{"js": "// The two paragraphs whose explicit \"12pt\" direct formatting (w:sz / w:szCs,\n// both worth 24 half-points = 12pt) must be dropped from BOTH the paragraph\n// mark's run properties (w:pPr/w:rPr) and the visible run's properties\n// (w:r/w:rPr) \u2014 while keeping everything else (fonts, text, paragraph\n// identity) untouched. We locate each paragraph by its distinctive text and\n// rewrite it in place via insertOoxml so we can surgically drop just the\n// two size tags instead of clobbering unrelated formatting.\nconst targets = [\n  \"Voici un bref aper\u00e7u de ce que la plupart des managers attendent d\\u2019un d\\u00e9veloppeur une fois que celui-ci est int\\u00e9gr\\u00e9 \\u00e0 son \\u00e9quipe.\",\n  \"On ne va pas se mentir, on embauche avant tout un d\\u00e9veloppeur pour qu\\u2019il puisse\\u2026 d\\u00e9velopper. Ou programmer, si l\\u2019on pr\\u00e9f\\u00e8re. Et pour cela, peu importe son niveau de comp\\u00e9tences de base, il est absolument indispensable d\\u2019avoir un implacable sens de la logique. \"\n];\n\nfor (const target of targets) {\n  // Search using a shortened, unique snippet so exact whitespace/trailing\n  // punctuation differences don't prevent a match.\n  const snippet = target.substring(0, 40);\n  const results = context.document.body.search(snippet, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target paragraph for: \" + snippet);\n  }\n\n  const paragraph = results.items[0].paragraphs.getFirst();\n  paragraph.load(\"text\");\n  await context.sync();\n\n  // Grab the paragraph's current OOXML so we can edit it directly and\n  // preserve every attribute we are not intentionally changing.\n  const ooxmlResult = paragraph.getOoxml();\n  await context.sync();\n\n  let xml = ooxmlResult.value;\n\n  // Pull just the <w:p \u2026>\u2026</w:p> fragment out of the pkg:package wrapper.\n  const pStart = xml.indexOf(\"<w:p \");\n  const pEnd = xml.indexOf(\"</w:p>\") + \"</w:p>\".length;\n  let pFragment = xml.substring(pStart, pEnd);\n\n  // Drop the explicit 12pt sizing from both the paragraph mark's rPr and\n  // the run's rPr, leaving rFonts (and everything else) intact.\n  pFragment = pFragment.split('<w:sz w:val=\"24\"/>').join(\"\");\n  pFragment = pFragment.split('<w:szCs w:val=\"24\"/>').join(\"\");\n\n  const wrappedOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + pFragment + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  paragraph.insertOoxml(wrappedOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Remove the explicit 12pt direct formatting (w:sz / w:szCs, both \"24\" in\n# half-points) from the two paragraphs below. The size tags must disappear\n# from BOTH the paragraph mark's run properties (w:pPr/w:rPr) and the\n# visible run's properties (w:r/w:rPr), while every other attribute\n# (fonts, text, paragraph identity, rsids, \u2026) stays exactly as-is. We find\n# each paragraph with Find.Execute on a short ASCII-safe snippet (accented\n# characters don't round-trip reliably through this console), fetch its\n# OOXML, strip just the two size tags, and write it back with InsertXML so\n# the edit is surgical instead of touching unrelated formatting.\n$d = $word.ActiveDocument\n\n$snippets = @(\n    \"Voici un bref aper\",\n    \"On ne va pas se mentir\"\n)\n\nforeach ($snippet in $snippets) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($snippet)\n    if (-not $found) {\n        throw \"Could not find target paragraph for snippet: $snippet\"\n    }\n\n    $para = $rng.Paragraphs(1)\n    $xml = $para.Range.WordOpenXML\n    $newXml = $xml.Replace('<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>', '')\n    $para.Range.InsertXML($newXml)\n}\n"}
